$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation")

$data = @{
    2  = @([double]"1.781018743393181E-15", [double]"4.057255051223404E-22", [double]"1.0700961573157958E-62")
    3  = @(0.03641868368847227, 19.12835270551491, 0.001239938123641662)
    4  = @(0.03814118378129024, 29.026306881176907, 33.7142861575064)
    5  = @(0.03778826035639545, 29.34969224937014, 43.912774468368156)
    6  = @(0.03789100862820734, 28.366955385060376, 46.337484667189464)
    7  = @(0.038998158757195554, 25.294465279055657, 48.238660452053104)
    8  = @(0.04291751810464542, 18.3269325394798, 49.1448849657594)
    9  = @(0.04574189809598686, 15.541395920050864, 49.11122521395621)
    10 = @(0.04769772290282887, 14.584029862084654, 48.98838892804203)
    11 = @(0.051140033815074114, 14.130623099305257, 48.435643821175276)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
}
